# Applies the "first completed version of the thesis" edits to the
# List of Tables front-matter page:
#   1. Retarget the TOC heading's own bookmark to a new Toc id.
#   2. Drop the stray "_GoBack" bookmark sitting inside "Table" and
#      re-merge the split "Ta"/"ble" runs into a single "Table" run.
#   3. Bump the cached PAGEREF page numbers for every table TOC entry.
#   4. Update the cached roman-numeral PAGE field result in the
#      section's primary footer from "i" to "x".

$d = $word.ActiveDocument

# 1. Rename bookmark _Toc439680270 -> _Toc440031416 (wraps "LIST OF TABLES").
$oldHeadingBookmark = $d.Bookmarks("_Toc439680270")
$headingRange = $oldHeadingBookmark.Range
$oldHeadingBookmark.Delete()
$d.Bookmarks.Add("_Toc440031416", $headingRange) | Out-Null

# 2. Remove the "_GoBack" bookmark and merge "Ta" + "ble" back into "Table".
$goBack = $d.Bookmarks("_GoBack")
$tableRange = $goBack.Range
$goBack.Delete()
# Expand slightly so Find can match the full word even though it currently
# lives in two separate runs ("Ta" / "ble").
$tableLineStart = $tableRange.Start - 2
$tableLineEnd = $tableRange.End + 3
$tableFindRange = $d.Range($tableLineStart, $tableLineEnd)
$tableFindRange.Find.Execute("Table", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Table", 2) | Out-Null

# 3. Update the cached PAGEREF results for each "Table" TOC line, in order.
$pageUpdates = @(
    "46",  "47";
    "53",  "54";
    "59",  "60";
    "62",  "63";
    "78",  "79";
    "82",  "83";
    "106", "93";
    "106", "108";
    "160", "161";
    "161", "162";
    "162", "163";
    "163", "164";
    "171", "172";
    "172", "173";
    "172", "173";
    "173", "174";
    "174", "175";
    "175", "176"
)

$fields = $d.Fields
$fieldIndex = 0
for ($i = 1; $i -le $fields.Count; $i++) {
    $field = $fields.Item($i)
    if ($field.Code.Text.Trim().StartsWith("PAGEREF")) {
        $oldVal = $pageUpdates[$fieldIndex * 2]
        $newVal = $pageUpdates[$fieldIndex * 2 + 1]
        $resultRange = $field.Result
        $resultRange.Find.Execute($oldVal, $false, $false, $false, $false, $false, `
            $true, 1, $false, $newVal, 2) | Out-Null
        $fieldIndex = $fieldIndex + 1
    }
}

# 4. Update the cached roman-numeral page number in the primary footer
#    of the first section from "i" to "x".
$footer = $d.Sections.Item(1).Footers.Item(1)
$pageField = $footer.Range.Fields.Item(1)
$pageResult = $pageField.Result
$pageResult.Find.Execute("i", $false, $false, $false, $false, $false, `
    $true, 1, $false, "x", 2) | Out-Null
